# Update crypto price/volume figures to the latest scrape snapshot
# (commit: "Updated symbol list on Fri Feb 10 15:39:25 UTC 2023 with GitHub Actions").
# Values are stored as text in the sheet (Price/Volume columns), so each
# assignment uses a leading apostrophe to force text entry and avoid Excel
# reinterpreting the numeric-looking / percentage-looking strings as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.47"
$ws.Range("E2").Value = "'-3.88%"
$ws.Range("D3").Value = "'40.38"
$ws.Range("E3").Value = "'-5.35%"
$ws.Range("D4").Value = "'5.083"
$ws.Range("D5").Value = "'0.07687"
$ws.Range("E5").Value = "'-6.12%"
$ws.Range("D6").Value = "'4.249"
$ws.Range("E6").Value = "'-2.68%"
$ws.Range("D7").Value = "'1.616"
$ws.Range("E7").Value = "'-8.71%"
$ws.Range("D8").Value = "'0.9170"
$ws.Range("E8").Value = "'-3.47%"
$ws.Range("D9").Value = "'0.1041"
$ws.Range("E9").Value = "'-6.90%"
$ws.Range("D10").Value = "'0.1771"
$ws.Range("E10").Value = "'-5.60%"
$ws.Range("D11").Value = "'0.09316"
$ws.Range("E11").Value = "'-0.41%"
$ws.Range("D12").Value = "'0.04437"
$ws.Range("E12").Value = "'-5.81%"
$ws.Range("D13").Value = "'0.1056"
$ws.Range("E13").Value = "'0.00%"
$ws.Range("D14").Value = "'0.001254"
$ws.Range("E14").Value = "'-2.24%"
$ws.Range("D15").Value = "'0.005814"
$ws.Range("E15").Value = "'1.77%"
$ws.Range("E16").Value = "'2,406.26%"
$ws.Range("E17").Value = "'0.28%"
$ws.Range("D18").Value = "'2.456"
$ws.Range("E18").Value = "'-3.00%"
$ws.Range("D19").Value = "'0.3318"
$ws.Range("E19").Value = "'-1.32%"
$ws.Range("D20").Value = "'6.943"
$ws.Range("E20").Value = "'-6.56%"
$ws.Range("E21").Value = "'0.69%"
$ws.Range("E22").Value = "'5.77%"
$ws.Range("E23").Value = "'-0.80%"
$ws.Range("D24").Value = "'0.001206"
$ws.Range("E24").Value = "'-3.00%"
$ws.Range("D25").Value = "'0.004109"
$ws.Range("E25").Value = "'-4.26%"
$ws.Range("D26").Value = "'0.0001301"
$ws.Range("E26").Value = "'6.26%"
$ws.Range("D38").Value = "'0.02481"
$ws.Range("E38").Value = "'-5.74%"
$ws.Range("D39").Value = "'0.05183"
$ws.Range("E39").Value = "'-7.81%"
$ws.Range("D40").Value = "'0.007945"
$ws.Range("E40").Value = "'-2.79%"
$ws.Range("E41").Value = "'-5.97%"
$ws.Range("D42").Value = "'0.007081"
$ws.Range("E42").Value = "'8.11%"
$ws.Range("D43").Value = "'0.001950"
$ws.Range("E43").Value = "'-7.90%"
$ws.Range("E44").Value = "'-2.64%"
$ws.Range("D45").Value = "'0.3068"
$ws.Range("E45").Value = "'-11.71%"
$ws.Range("D46").Value = "'0.00006433"
$ws.Range("E46").Value = "'-5.31%"
$ws.Range("E47").Value = "'-0.26%"
$ws.Range("D48").Value = "'0.003002"
$ws.Range("E48").Value = "'-27.04%"
$ws.Range("D49").Value = "'0.004558"
$ws.Range("E49").Value = "'36.15%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.26%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.26%"
